# Add a new server entry ("Server6") to the server/firewall inventory table,
# labelled with its node name and IP address, connecting to Server3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (vertical-center + wrap-text style) from the row above
# so the new row A10:B10 matches the rest of the Name/IP columns.
$ws.Range("A9:B9").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A10").Value = "Server6"
$ws.Range("B10").Value = "192.168.8.1"
$ws.Range("C10").Value = "Server3"

# Match the author's final selection.
$ws.Range("C13").Select()
